# Insert two new weekly price rows above existing row 246, pushing the
# existing rows 246-249 down to 248-251 (their contents are preserved as-is).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 246 (this shifts row 246 and below down by 2).
$ws.Rows.Item(246).Insert()
$ws.Rows.Item(246).Insert()

# New row 246: Coliflor, Primera, week of 2022-04-05 (serial 44656)
$ws.Range("A246").Value = 7
$ws.Range("B246").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C246").Value = "Ñuble"
$ws.Range("D246").Value = 44656
$ws.Range("E246").Value = 16
$ws.Range("F246").Value = 100112008
$ws.Range("G246").Value = "Coliflor"
$ws.Range("H246").Value = "Sin especificar"
$ws.Range("I246").Value = "Primera"
$ws.Range("J246").Value = 120
$ws.Range("K246").Value = 1200
$ws.Range("L246").Value = 1300
$ws.Range("M246").Value = 1250
$ws.Range("N246").Value = "$/unidad"
$ws.Range("O246").Value = "Región del Maule"
$ws.Range("P246").Value = 1250
$ws.Range("Q246").Value = 1
$ws.Range("R246").Value = "Hortaliza"

# New row 247: Coliflor, Segunda, week of 2022-04-05 (serial 44656)
$ws.Range("A247").Value = 7
$ws.Range("B247").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C247").Value = "Ñuble"
$ws.Range("D247").Value = 44656
$ws.Range("E247").Value = 16
$ws.Range("F247").Value = 100112008
$ws.Range("G247").Value = "Coliflor"
$ws.Range("H247").Value = "Sin especificar"
$ws.Range("I247").Value = "Segunda"
$ws.Range("J247").Value = 60
$ws.Range("K247").Value = 1000
$ws.Range("L247").Value = 1000
$ws.Range("M247").Value = 1000
$ws.Range("N247").Value = "$/unidad"
$ws.Range("O247").Value = "Región del Maule"
$ws.Range("P247").Value = 1000
$ws.Range("Q247").Value = 1
$ws.Range("R247").Value = "Hortaliza"

# Ensure the date cells keep the workbook's date number format (style index 2
# used throughout column D), matching the formatting of the rest of the column.
$ws.Range("D246").NumberFormat = $ws.Range("D248").NumberFormat
$ws.Range("D247").NumberFormat = $ws.Range("D248").NumberFormat
